$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '64.879.96', '  -2.99%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.440.82', '  -2.80%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.16%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '571.23', '  +0.53%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '175.21', '  -7.05%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.623', '  +0.52%  '),
    @(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.06%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.624', '  -1.68%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.159', '  +4.89%  '),
    @(11, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '54.91', '  +0.26%  '),
    @(12, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000273', '  +0.98%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '9.11', '  -3.44%  '),
    @(14, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.990.79', '  -2.54%  '),
    @(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.121', '  -0.88%  '),
    @(16, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.448.46', '  -2.47%  '),
    @(17, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '18.08', '  -1.12%  '),
    @(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '11.86', '  -1.55%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '64.884.61', '  -2.95%  '),
    @(20, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.987', '  -1.14%  '),
    @(21, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '407.26', '  -4.83%  '),
    @(22, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '4.17', '  -0.52%  '),
    @(23, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '4.45', '  +7.57%  '),
    @(24, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '13.47', '  +9.56%  '),
    @(25, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '83.61', '  -1.91%  '),
    @(26, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '10.80', '  -3.03%  '),
    @(27, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.80', '  -3.43%  '),
    @(28, 'LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '6.00', '  -2.34%  '),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '8.98', '  -2.89%  '),
    @(30, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '29.89', '  -1.93%  '),
    @(31, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.58', '  -0.47%  '),
    @(32, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.53', '  -1.99%  '),
    @(33, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '584.53', '  -9.14%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.108', '  -3.62%  '),
    @(35, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '59.72', '  -0.30%  '),
    @(36, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.152', '  +2.53%  '),
    @(37, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.02%  '),
    @(38, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '3.57', '  +5.88%  '),
    @(39, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0774', '  -4.80%  '),
    @(40, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '36.20', '  -6.22%  '),
    @(41, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.376', '  -4.11%  '),
    @(42, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '3.170.68', '  +4.35%  '),
    @(43, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.999', '  +0.03%  '),
    @(44, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.93', '  +1.68%  '),
    @(45, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '2.50', '  -6.49%  '),
    @(46, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '3.22', '  -3.84%  '),
    @(47, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0410', '  -2.39%  '),
    @(48, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.131', '  -1.20%  '),
    @(49, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.63', '  -4.68%  '),
    @(50, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '8.45', '  -2.17%  '),
    @(51, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '136.85', '  -3.53%  '),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Price column (D) often looks numeric (e.g. '1.00', '0.623'); force it to
    # stay a text value like the original inlineStr cell, then strip the
    # temporary text NumberFormat so no stray style index is left behind.
    $dcell = $ws.Cells.Item($r, 4)
    $dcell.NumberFormat = '@'
    $dcell.Value = $row[3]
    $dcell.Style = 'Normal'

    $ws.Cells.Item($r, 5).Value = $row[4]
}